# New crime data collected — weekly CompStat report refresh
# (Volume/week-range header text + the 120th Precinct crime-complaint table)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume number + reporting week range ---
$ws.Range("A8").Value = "Volume 30   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/25/2023  Through  12/31/2023"

# --- Numeric -> Text placeholder conversions ("0" or "***.*") ---
# A leading apostrophe forces Excel to store the value as literal text
# instead of re-parsing it as a number; PasteSpecial(formats) then copies
# the donor cell's number format/style (General, right/center) onto the
# target so the cell ends up identical to its text-placeholder siblings.
function Set-TextPlaceholder($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range("C14").Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
}

Set-TextPlaceholder "G14" "0"
Set-TextPlaceholder "H14" "***.*"

Set-TextPlaceholder "C15" "0"

Set-TextPlaceholder "C27" "0"
Set-TextPlaceholder "D27" "0"
Set-TextPlaceholder "E27" "***.*"

Set-TextPlaceholder "G28" "0"
Set-TextPlaceholder "H28" "***.*"

Set-TextPlaceholder "G29" "0"
Set-TextPlaceholder "H29" "***.*"

# --- Plain numeric value updates ---
$ws.Range("M14").Value = -7.692307692307

$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("M15").Value = -54.761904761904

$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 10
$ws.Range("I16").Value = 163
$ws.Range("J16").Value = 117
$ws.Range("K16").Value = 39.316239316239
$ws.Range("L16").Value = 89.534883720930
$ws.Range("M16").Value = -42.402826855123
$ws.Range("N16").Value = -83.956692913385

$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 66.666666666666
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = 47.826086956521
$ws.Range("I17").Value = 470
$ws.Range("J17").Value = 394
$ws.Range("K17").Value = 19.289340101522
$ws.Range("L17").Value = 70.289855072463
$ws.Range("M17").Value = 59.322033898305
$ws.Range("N17").Value = -40.581542351453

$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -17.647058823529
$ws.Range("I18").Value = 142
$ws.Range("J18").Value = 124
$ws.Range("K18").Value = 14.516129032258
$ws.Range("L18").Value = 125.396825396825
$ws.Range("M18").Value = -54.487179487179
$ws.Range("N18").Value = -91.671554252199

$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -27.272727272727
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 33
$ws.Range("H19").Value = -15.151515151515
$ws.Range("I19").Value = 413
$ws.Range("J19").Value = 347
$ws.Range("K19").Value = 19.020172910662
$ws.Range("L19").Value = 41.438356164383
$ws.Range("M19").Value = -0.481927710843
$ws.Range("N19").Value = -29.881154499151

$ws.Range("C20").Value = 1
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 156
$ws.Range("J20").Value = 136
$ws.Range("K20").Value = 14.705882352941
$ws.Range("L20").Value = 102.597402597403
$ws.Range("M20").Value = -18.75
$ws.Range("N20").Value = -88.418708240534

$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 101
$ws.Range("G21").Value = 99
$ws.Range("H21").Value = 2.020202020202
$ws.Range("I21").Value = 1375
$ws.Range("J21").Value = 1139
$ws.Range("K21").Value = 20.71992976295
$ws.Range("L21").Value = 66.868932038835
$ws.Range("M21").Value = -11.404639175257
$ws.Range("N21").Value = -75.113122171945

$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -25
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 86
$ws.Range("J23").Value = 58
$ws.Range("K23").Value = 48.275862068965
$ws.Range("L23").Value = 65.384615384615
$ws.Range("M23").Value = 68.627450980392

$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 4.166666666666
$ws.Range("G24").Value = 114
$ws.Range("H24").Value = -5.263157894736
$ws.Range("I24").Value = 1315
$ws.Range("J24").Value = 1206
$ws.Range("K24").Value = 9.038142620232
$ws.Range("L24").Value = 48.085585585585
$ws.Range("M24").Value = -12.391738840772

$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 44
$ws.Range("H25").Value = -9.090909090909
$ws.Range("I25").Value = 724
$ws.Range("J25").Value = 636
$ws.Range("K25").Value = 13.836477987421
$ws.Range("L25").Value = 31.636363636363
$ws.Range("M25").Value = -41.233766233766

$ws.Range("C26").Value = 2
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 66.666666666666
$ws.Range("I26").Value = 31
$ws.Range("K26").Value = -8.823529411764
$ws.Range("L26").Value = -18.421052631578

$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 400
$ws.Range("L27").Value = 48.387096774193

$ws.Range("M28").Value = -17.857142857142
$ws.Range("N28").Value = -79.464285714285

$ws.Range("M29").Value = -16
$ws.Range("N29").Value = -79
